# Applies the "Core Functionality section finished" edit:
#  - Turns the blank title-page placeholder paragraph into a new
#    "Introduction" Heading 1 (followed by 5 blank spacer paragraphs),
#    mirroring the blank-spacer pattern used by the other H1 sections.
#  - Restyles "1 Project Setup" (drops the stray lastRenderedPageBreak,
#    adds the white/background1 color run property to its pPr).
#  - Fills in the body copy for section 2 (Core Functionality) and its
#    three subsections (2.1 Implementation Features, 2.2 Choice of
#    Algorithms, 2.3 Modular Components), removing the blank spacer
#    paragraphs that used to be placeholders for this text.
#  - Moves the lastRenderedPageBreak marker from "3 Documentation /
#    Implementation of Code" onto "2.1 Implementation Features", since
#    that's where the page now actually breaks once the new body text
#    is in place.

$d = $word.ActiveDocument

function Find-ParaIndex([string]$text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*$text*") {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# Work from the bottom of the document upwards so that each edit's
# paragraph-index lookups aren't invalidated by earlier insertions.
# ------------------------------------------------------------------

# 8. "3 Documentation / Implementation of Code" heading: drop lastRenderedPageBreak
$idx = Find-ParaIndex("3 Documentation / Implementation of Code")
$d.Paragraphs.Item($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
      </w:pPr>
      <w:r>
        <w:t>3 Documentation / Implementation of Code</w:t>
      </w:r>
    </w:p>')

# 7. "2.3 Modular Components": replace the 3 blank paragraphs after it with the new body text
$idx = Find-ParaIndex("2.3 Modular Components")
$p1 = $d.Paragraphs.Item($idx + 1)
$p2 = $d.Paragraphs.Item($idx + 3)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$rng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">The Phase 1 implementation follows a modular design to encourage simplicity and promote code reuse. Core gameplay entities are encapsulated within individual classes, including Palle, Ball, Block and Level. Each class is responsible for managing its own state, </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>update</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> logic, and rendering behavior.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>The Level class serves as a structural coordinator rather than a rendering or physics engine. It encapsulates level specific data such as block layout and type. Interfaces are provided for updating and querying current level state.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>Game flow is managed by a central game loop which coordinates input handling, state updates, collision checks, and rendering. Game states such as active gameplay and Game Over are represented by distinct logical components, during Phase 2 clean transitions will be enabled by moving and expanding upon the state code within the manger system.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>Supporting systems, such as scoring and life tracking, are implemented within the core game loop and during Phase 2 will be moved into their own managers to avoid tight coupling with gameplay entities. Configuration values, including movement speeds, screen dimensions, and initial lives, are centralized to simplify tuning and maintenance.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>This modular architecture improves readability, facilitates debugging, and provides a flexible foundation to simplify future phases of development.</w:t>
      </w:r>
    </w:p>')

# 6. "2.2 Choice of Algorithms": replace the blank paragraph after it with the new body text
$idx = Find-ParaIndex("2.2 Choice of Algorithms")
$p1 = $d.Paragraphs.Item($idx + 1)
$p1.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>The algorithms selected for Phase 1 prioritize simplicity, determinism, and clarity to support both gameplay consistency and maintainability. Paddle movement is calculated using linear translation based upon player input and a configurable speed value and clamped to the horizontal bounds of the screen.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>Ball movement is handled using basic vector arithmetic, where position updates are computed as a function of velocity and delta time. Collision detection uses axis-aligned bounding box (AABB) comparisons to determine intersections between the ball and other game objects such as blocks and the paddle. Collision response is implemented by inverting the appropriate component of the velocity vector depending on the collision surface.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Block management within a level is supported by a simple collection-based algorithm. The Level class maintains a collection of active blocks and updates the collection as blocks are destroyed. Level completion is determined by checking whether all destructible blocks have been destroyed, triggering either level completion logic or game termination depending on the current scope.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>Scoring and life management use simple counter-based algorithms. Block destruction triggers immediate score updates, while life loss events decrement a life counter and reset the ball state. These algorithms are computationally lightweight and well suited to a real-time 2D arcade game.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>Use of the chosen algorithms reflects the team’s desire to create a faithful reconstruction of the classic Breakout arcade game. Each algorithm is similar, if not the same, as the original game and ensures the game will behave as the user expects it should.</w:t>
      </w:r>
    </w:p>')

# 5. "2.1 Implementation Features": replace the 2 blank paragraphs after it with the new body text
$idx = Find-ParaIndex("2.1 Implementation Features")
$p1 = $d.Paragraphs.Item($idx + 1)
$p2 = $d.Paragraphs.Item($idx + 2)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$rng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>Phase 1 implementation is focused on foundational gameplay mechanics and systems required for a playable experience. The paddle movement allows the player to control a horizontal paddle using keyboard input. Movement is smooth and frame-rate independent, ensuring consistent play experience.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>Ball movement is implemented using a two-dimensional velocity vector updated each frame based on elapsed time. The ball responds to collisions with the paddle, bricks, and screen boundaries through deterministic reflection logic. When the ball falls below the bottom void boundary, a life is deducted and the ball is reset to the starting position.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>Block collision functionality supports a basic block type that is destroyed upon a single collision with the ball. When a block is destroyed, the scoring system is updated. The HUD provides a scoring counter and indicator with the player’s remaining lives</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> during game play. Game state transitions, including a new game and triggering the Game Over state when all lives are lost are not implemented as part of the immediate game play loop. State changes will be implemented to allow for game over, start game, and menu states within the next phase.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>A dedicated Level class is implemented to manage level structure and progression. The Level class is responsible for loading block layouts, which are provided in JSON. The Level class will instantiate blocks, track remaining blocks, and determine when a level has been completed. Abstracting in this way allows level definitions to be modified or extended independently of core gameplay logic. For Phase 1, only a single playable level is implemented.</w:t>
      </w:r>
    </w:p><w:p>
      <w:r>
        <w:t>Combined, these features establish a complete and playable core which will serve as the baseline game for future enhancements. During the next phase the team will be adding state transitions, further modularizing the code, and including menus which will enable the addition of options and configuration modifications.</w:t>
      </w:r>
    </w:p>')

# 4. "2.1 Implementation Features" heading: add lastRenderedPageBreak
$idx = Find-ParaIndex("2.1 Implementation Features")
$d.Paragraphs.Item($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>2.1 Implementation Features</w:t>
      </w:r>
    </w:p>')

# 3. "2 Core Functionality": replace the first blank paragraph after it with the section intro text
$idx = Find-ParaIndex("2 Core Functionality")
$p1 = $d.Paragraphs.Item($idx + 1)
$p1.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>This section describes the core functional elements of the Breakout game implementation during Phase 1 of the project. Phase 1 focuses on establishing a stable and playable foundation of the Breakout-style game by implementing essential gameplay mechanics, algorithms, and modular components. The goal of this phase is to deliver a functional game phase of the MVP</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> demonstrating a playable game, maintainable structure and extensibility to build on for future phases.</w:t>
      </w:r>
    </w:p>')

# 2. "1 Project Setup" heading: drop lastRenderedPageBreak, add white/background1 color
$idx = Find-ParaIndex("1 Project Setup")
$d.Paragraphs.Item($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:rPr>
          <w:color w:val="FFFFFF" w:themeColor="background1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>1 Project Setup</w:t>
      </w:r>
    </w:p>')

# 1. Blank title-page placeholder paragraph (just before "1 Project Setup"):
#    turn it into the new "Introduction" Heading 1, then add 5 blank spacer
#    paragraphs after it (matching the other H1 sections' spacing pattern).
$idx = Find-ParaIndex("1 Project Setup")
$introIdx = $idx - 1
$introP = $d.Paragraphs.Item($introIdx)
$introP.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Introduction</w:t>
      </w:r>
    </w:p>')
$introP2 = $d.Paragraphs.Item($introIdx)
$introP2.Range.InsertParagraphAfter()
$afterIntro = $d.Paragraphs.Item($introIdx + 1)
$afterIntro.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p/><w:p/><w:p/><w:p/>')

Write-Output "DONE"
